$wb = $excel.ActiveWorkbook

# --- "About" section: OrangeHRM -> OrangeHRM OS 5.7 -------------------------
$profileSheet = $wb.Worksheets.Item("Profile")
$profileSheet.Range("B2").Value = "OrangeHRM OS 5.7"

# --- Login sheet: remove the login-URL hyperlink and its one-off font -------
$loginSheet = $wb.Worksheets.Item("Login")
$linkCell = $loginSheet.Range("A2")

# Drop the hyperlink itself (and the relationship it points at).
$loginSheet.Hyperlinks.Delete()

# The hyperlink also carried a dedicated underlined/blue font; restore the
# plain font used by the rest of the sheet now that it is no longer a link.
$loginSheet.Range("B1").Copy()
$linkCell.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
